$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 4 (pushes existing rows 4..43 -> 5..44)
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly price record
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C4").Value = 'Arica y Parinacota'
$ws.Range("D4").Value = 44537
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 100112031
$ws.Range("G4").Value = 'Poroto verde'
$ws.Range("H4").Value = 'Sin especificar'
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 1700
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = 550
$ws.Range("N4").Value = '$/kilo'
$ws.Range("O4").Value = 'Región de Arica y Parinacota'
$ws.Range("P4").Value = 550
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 'Hortaliza'
